$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (row => @{ D=...; J=...; K=...; L=...; M=...; P=... })
$rowData = @{
    3  = @{ D = 44379; J = 120; K = 12000; L = 13000; M = 12667; P = 974 }
    4  = @{ D = 44580; J = 160; K = 11000; L = 12000; M = 11500; P = 885 }
    5  = @{ D = 44469; J = 140; K = 13000; L = 14000; M = 13500; P = 1038 }
    6  = @{ D = 44406; J = 160; K = 17000; L = 18000; M = 17500; P = 1346 }
    7  = @{ D = 44592; J = 120; K = 12000; L = 13000; M = 12500; P = 962 }
    8  = @{ D = 44159; J = 100; K = 23000; L = 24000; M = 23500; P = 1808 }
    9  = @{ D = 44320; J = 160; K = 19000; L = 20000; M = 19500; P = 1500 }
    10 = @{ D = 44397; J = 140; K = 12500; L = 13000; M = 12750; P = 981 }
    11 = @{ D = 44229; J = 120; K = 44000; L = 45000; M = 44500; P = 3423 }
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("J$r").Value = $vals.J
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("P$r").Value = $vals.P
}
